# OPAR Ordnance.xlsx - "CURRENT AS OF" rolled from D3.1 to D3.2; new expenditures
# entered for that day, and the just-consumed G19:Q19 ordnance block was
# re-centered (matches clicking the "Center" alignment button over that range).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- advance the "CURRENT AS OF" marker from D3.1 to D3.2 -------------------
$ws.Range("R1").Value = "D3.2"

# --- new expenditures logged against the D3.1 / D3.2 columns ---------------
$ws.Range("H4").Value  = 3   # AGM-84D (Harpoon)  - D3.1
$ws.Range("I12").Value = 8   # GBU-31             - D3.2
$ws.Range("I14").Value = 8   # GBU-10             - D3.2
$ws.Range("I18").Value = 8   # GBU-31(V) 3/B      - D3.2
$ws.Range("I19").Value = 2   # Fuel tanks         - D3.2

# --- re-center the D3.1:D7.2 block of the "Fuel tanks" row ------------------
# L19 carried a left/indent alignment (from its font/merge heritage) - drop
# the indent before the block-wide center so it ends up plain "center" like
# its neighbours.
$ws.Range("L19").IndentLevel = 0
$ws.Range("G19:Q19").HorizontalAlignment = -4108   # xlCenter

# --- move the active selection ----------------------------------------------
$ws.Range("I17").Select()
